$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6159312725067139
$ws.Range("B1").Value = 1.580912709236145
$ws.Range("C1").Value = 3.874247312545776
$ws.Range("D1").Value = 5.390730381011963
$ws.Range("E1").Value = 1.609069466590881
